# Update gh-pages to output generated at 456a3b4
# Refresh the scraped event-listing data on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets: the oldest event (2024-04-11) has aged
# out of the feed, every remaining row's stats were re-scraped, and a new
# upcoming entry shifts into view.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)

# Row 2: 南宁·2024三月三国潮动漫节（良牙春典）
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2024-05-01"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Value = "南宁·2024三月三国潮动漫节（良牙春典）"
$ws.Cells.Item(2, 4).Value = "民族大道106号 南宁国际会展中心"
$ws.Cells.Item(2, 5).Value = "2024.05.01 09:30-05.02 17:30"
$ws.Cells.Item(2, 6).Value = 4622
$ws.Cells.Item(2, 7).Value = 55
$ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82416"
$ws.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"

# Row 3: 南宁·原x穹x崩only
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2024-05-19"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 3).Value = "南宁·原x穹x崩only"
$ws.Cells.Item(3, 4).Value = "明秀东路157号 利泰国际大酒店"
$ws.Cells.Item(3, 5).Value = "2024.05.19 10:00-05.19 17:00"
$ws.Cells.Item(3, 6).Value = 136
$ws.Cells.Item(3, 7).Value = 35
$ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83070"
$ws.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

# Row 4: 南宁·布谷鸟动漫展4th
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2024-06-09"
$ws.Cells.Item(4, 2).Style = "Normal"
$ws.Cells.Item(4, 3).Value = "南宁·布谷鸟动漫展4th"
$ws.Cells.Item(4, 4).Value = "亭洪路45号 百益上河城"
$ws.Cells.Item(4, 5).Value = "2024.06.09 10:00-06.10 17:00"
$ws.Cells.Item(4, 6).Value = 805
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82241"
$ws.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

# The oldest entry (2024-04-11) has aged off the list; drop the now-obsolete trailing row.
$ws.Rows.Item(5).Delete()

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)

# Row 2: 南宁·2024三月三国潮动漫节（良牙春典）
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2024-05-01"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Value = "南宁·2024三月三国潮动漫节（良牙春典）"
$ws.Cells.Item(2, 4).Value = "民族大道106号 南宁国际会展中心"
$ws.Cells.Item(2, 5).Value = "2024.05.01 09:30-05.02 17:30"
$ws.Cells.Item(2, 6).Value = 4622
$ws.Cells.Item(2, 7).Value = 55
$ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82416"
$ws.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"

# Row 3: 南宁·原x穹x崩only
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2024-05-19"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 3).Value = "南宁·原x穹x崩only"
$ws.Cells.Item(3, 4).Value = "明秀东路157号 利泰国际大酒店"
$ws.Cells.Item(3, 5).Value = "2024.05.19 10:00-05.19 17:00"
$ws.Cells.Item(3, 6).Value = 136
$ws.Cells.Item(3, 7).Value = 35
$ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83070"
$ws.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

# Row 4: 南宁·布谷鸟动漫展4th
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2024-06-09"
$ws.Cells.Item(4, 2).Style = "Normal"
$ws.Cells.Item(4, 3).Value = "南宁·布谷鸟动漫展4th"
$ws.Cells.Item(4, 4).Value = "亭洪路45号 百益上河城"
$ws.Cells.Item(4, 5).Value = "2024.06.09 10:00-06.10 17:00"
$ws.Cells.Item(4, 6).Value = 805
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82241"
$ws.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

# Row 5: 南宁·浪漫古典·百年经典世界名曲音乐会
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "2024-06-22"
$ws.Cells.Item(5, 2).Style = "Normal"
$ws.Cells.Item(5, 3).Value = "南宁·浪漫古典·百年经典世界名曲音乐会"
$ws.Cells.Item(5, 4).Value = "广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅"
$ws.Cells.Item(5, 5).Value = "2024.06.22 20:00-06.22 21:30"
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 50
$ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83959"
$ws.Cells.Item(5, 9).Value = "//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg"

# The oldest entry (2024-04-11) has aged off the list; drop the now-obsolete trailing row.
$ws.Rows.Item(6).Delete()
